$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "2007年" row). This shifts all subsequent rows
# (2010年, 2012年, 2015年, 2017年) up by one, matching the target layout.
$ws.Rows.Item(2).Delete()
